# Update "想去人数" (want-to-go count) figures on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet (column F = 想去人数)
$wsExhibition.Range("F4").Value  = 113
$wsExhibition.Range("F7").Value  = 972
$wsExhibition.Range("F9").Value  = 2063
$wsExhibition.Range("F10").Value = 1074
$wsExhibition.Range("F14").Value = 367
$wsExhibition.Range("F19").Value = 1520
$wsExhibition.Range("F21").Value = 676
$wsExhibition.Range("F22").Value = 562
$wsExhibition.Range("F23").Value = 12021
$wsExhibition.Range("F24").Value = 12031
$wsExhibition.Range("F26").Value = 678
$wsExhibition.Range("F29").Value = 299

# 全部类型 sheet (column F = 想去人数)
$wsAllTypes.Range("F6").Value  = 113
$wsAllTypes.Range("F9").Value  = 972
$wsAllTypes.Range("F11").Value = 2063
$wsAllTypes.Range("F12").Value = 1074
$wsAllTypes.Range("F16").Value = 367
$wsAllTypes.Range("F23").Value = 1520
$wsAllTypes.Range("F25").Value = 676
$wsAllTypes.Range("F26").Value = 562
$wsAllTypes.Range("F27").Value = 12021
$wsAllTypes.Range("F28").Value = 12031
$wsAllTypes.Range("F30").Value = 678
$wsAllTypes.Range("F33").Value = 299
